$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.054.27"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.550.14"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.04"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "197.52"
$ws.Range("E6").Value = "  +6.41%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.654"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.17"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000304"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.57"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.118.77"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "603.67"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.25"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.230.14"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.72"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.548.82"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.10"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.31"
$ws.Range("E23").Value = "  +6.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.87"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.67"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.89"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.45"
$ws.Range("E30").Value = "  +22.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.17"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.70"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.40"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0841"
$ws.Range("E35").Value = "  +8.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.789.92"
$ws.Range("E36").Value = "  +7.40%  "
$ws.Range("E37").Value = "  -4.59%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.70"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.396"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.81"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "488.91"
$ws.Range("E42").Value = "  -7.80%  "
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("E45").Value = "  -3.84%  "
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  -4.15%  "
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.23"
$ws.Range("E51").Value = "  -2.08%  "
